# Add the new "Creating the TCP Server for the Tournament" section at the
# end of the document, right after the existing "Creating the HTTP Server
# and the router" section (commit: "Added the very basics for tcp").

$d = $word.ActiveDocument

$headingText = "Creating the TCP Server for the Tournament (3h)"
$bodyText = "Since Im creating the TCP server ahead of the services which are supposed to send the tournament data I won" + [char]8217 + "t be able to properly test it. I have to admit that this was not a wise choice, but since I already started, I will try to create the basic structure and do the actually implementation in the tournament branch. Suddenly stopping and continuing with the services would probably cause thing like conflicts which would take me time to fix it, which I most definitely not have."

# Locate the end of the last paragraph of the "HTTP server" section (the
# one that ends "...in order to compare the routes efficiently." followed
# by a single trailing space run), then move past that trailing space so
# the new paragraphs are appended after it, not inside it.
$rng = $d.Content
$found = $rng.Find.Execute("in order to compare the routes efficiently.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor paragraph not found"
}
[void]$rng.MoveEnd(1, 1)
$rng.Collapse(0)

$rng.InsertAfter("`r" + $headingText + "`r" + $bodyText)

# The new heading paragraph currently has default (Normal) formatting;
# give it the "Heading 2" style used by the other section headings.
$d2 = $word.ActiveDocument
$count = $d2.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d2.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13)
    if ($text -eq $headingText) {
        $p.Style = "Heading 2"
        break
    }
}
